$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Helper: replace the paragraph that contains $findText with the OOXML given in
# $newParaXml (a full "<w:p ...>...</w:p>" fragment, including xmlns:w). This
# runtime's Range.InsertXML merges a paragraph-worth of OOXML inserted at a
# paragraph's start position into that paragraph (re-using its existing
# <w:pPr>), duplicating the paragraph's old runs after the newly inserted
# ones - so we follow up by deleting the now-duplicated original text that
# trails the freshly inserted content.
function Replace-Paragraph($findText, $newParaXml) {
    $rng = $d.Content
    $null = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $rng.Paragraphs(1)
    $start = $para.Range.Start
    $oldLen = $para.Range.End - $para.Range.Start - 1

    $ip = $d.Range($start, $start)
    $ip.InsertXML($newParaXml)

    $para2 = $d.Range($start, $start).Paragraphs(1)
    $newLen = $para2.Range.End - $para2.Range.Start - 1 - $oldLen

    $dupStart = $start + $newLen
    $dupEnd = $dupStart + $oldLen
    $d.Range($dupStart, $dupEnd).Delete()
}

# ---------------------------------------------------------------------------
# 1. New TODO item: "Update Call to use builtin AsmJit functionality for
#    remote code gen." right after "Forward declaration header.", taking
#    over the _GoBack bookmark that sat at the end of that paragraph.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$rng = $d.Content
$null = $rng.Find.Execute("Forward declaration header.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$ip = $d.Range($rng.End, $rng.End)
$newParaXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Update Call to use builtin AsmJit functionality for remote code gen.</w:t></w:r><w:bookmarkStart w:id=`"500`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"500`"/></w:p>"
$ip.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 2. Move the lastRenderedPageBreak marker from "Improved relative
#    instruction rebuilding ..." up to "Transactional hooking." (the extra
#    TODO line above pushed the rendered page break earlier).
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:lastRenderedPageBreak/><w:t>Transactional hooking.</w:t></w:r></w:p>"
Replace-Paragraph "Transactional hooking." $xml

$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:t>Improved relative instruction rebuilding (including conditionals).</w:t></w:r><w:r w:rsidR=`"00BD30BB`" w:rsidRPr=`"008B06FC`"><w:t xml:space=`"preserve`"> x64 has far more IP relative instructions than x86.</w:t></w:r></w:p>"
Replace-Paragraph "Improved relative instruction rebuilding" $xml

# ---------------------------------------------------------------------------
# 3. Move the lastRenderedPageBreak marker from "Improve export forwarding
#    code ..." up to "Full support for writing back to PE file ...".
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:lastRenderedPageBreak/><w:t>Full support for writing back to PE file, including automatically performing adjustments where required to fit in new data or remove unnecessary space.</w:t></w:r></w:p>"
Replace-Paragraph "Full support for writing back to PE file" $xml

$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:t>Improve export forwarding code to detect and handle forward-by-</w:t></w:r><w:r w:rsidR=`"00197559`" w:rsidRPr=`"008B06FC`"><w:t>ordinal</w:t></w:r><w:r w:rsidRPr=`"008B06FC`"><w:t xml:space=`"preserve`"> explicitly rather than forcing the user to detect it and do string manipulation and conversion. </w:t></w:r></w:p>"
Replace-Paragraph "Improve export forwarding code to detect" $xml

Write-Output "Done"
